$d = $word.ActiveDocument

# --- Locate the "Bugs:" paragraph and the review paragraph right after it ---
$bugsPara = $null
$reviewPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "Bugs:") {
        $bugsPara = $p
    } elseif (($bugsPara -ne $null) -and ($reviewPara -eq $null) -and ($t -like "My search (dog bones)*")) {
        $reviewPara = $p
    }
}

# --- Remove the stray _GoBack bookmark currently sitting on the "Bugs:" paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Find the split point right after "(??) " so the bookmark end lands there ---
$findRng = $d.Content
$findRng.Find.Execute("(??) that", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $findRng.Start + 5

# --- Re-anchor _GoBack so it wraps the review text from its start through "(??) " ---
$reviewStart = $reviewPara.Range.Start
$bmRange = $d.Range($reviewStart, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Append the trailing " FIXED" (in red) to the review paragraph ---
$reviewPara.Range.InsertAfter(" ")
$reviewPara.Range.InsertAfter("FIXED")

$pEnd = $reviewPara.Range.End
$fixedRange = $d.Range($pEnd - 1 - 5, $pEnd - 1)
$fixedRange.Font.Color = 255
